# Computer Assignment-SA 1.docx
# "Population proportion and other questions"
#
# The paragraph asking about SSSBE/SSI special incentives needs to be split
# into multiple runs:
#   1. " units, be it SSSBE or SSI, "                       (unchanged)
#   2. "that if the population proportion is less than 25%,"  (yellow highlight)
#   3. " there is a need for providing "                    (unchanged)
#   4. an empty "_GoBack" bookmark inserted here
#   5. "special incentives. ... or both?"                   (unchanged)

$d = $word.ActiveDocument

# --- Step 1: highlight "that if the population proportion is less than 25%," ---
$highlightRange = $d.Content
$highlightRange.Find.Execute(
    "that if the population proportion is less than 25%,",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null
$highlightRange.Font.HighlightColorIndex = 7   # wdYellow

# --- Step 2: insert the (empty) "_GoBack" bookmark right before "special incentives." ---
$gobackRange = $d.Content
$gobackRange.Find.Execute(
    "special incentives.  Based on your sample",
    $true, $false, $false, $false, $false, $true, 1, $false, "", 0
) | Out-Null
$gobackRange.Collapse(1)   # wdCollapseStart
$d.Bookmarks.Add("_GoBack", $gobackRange) | Out-Null
